$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 08:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6048317
$ws.Range("C4").Value = 1683
$ws.Range("D4").Value = 3348377
$ws.Range("E4").Value = 2515137
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 184803

# Row 6 - India
$ws.Range("B6").Value = 3392295
$ws.Range("C6").Value = 7720
$ws.Range("D6").Value = 2585030
$ws.Range("E6").Value = 745540
$ws.Range("G6").Value = 31
$ws.Range("H6").Value = 61725

# Row 57 - Kirguistan
$ws.Range("B57").Value = 43587
$ws.Range("C57").Value = 128
$ws.Range("D57").Value = 37726
$ws.Range("E57").Value = 4804

# Row 62 - Uzbekistan
$ws.Range("B62").Value = 40613
$ws.Range("C62").Value = 166
$ws.Range("E62").Value = 3113
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 300

# Row 72 - Australia
$ws.Range("D72").Value = 20633
$ws.Range("E72").Value = 4232

# Row 74 - Chequia
$ws.Range("B74").Value = 23300
$ws.Range("C74").Value = 131
$ws.Range("E74").Value = 5692

# Row 111 - Hong Kong
$ws.Range("E111").Value = 473
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 83

# Row 152 - Georgia
$ws.Range("B152").Value = 1455
$ws.Range("C152").Value = 8
$ws.Range("D152").Value = 1196
$ws.Range("E152").Value = 240

$wb.Save()
